# Generate Report for Handoff
# Updates the localization-status workbook to reflect that b.md has been
# handed off again: status moves from "Handed back: in sync with en-US" to
# "Ready for handoff", and a fresh handoff file/timestamp is recorded for
# both the zh-cn and de-de target sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$status = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet: row 3 is "b.md" -- refresh its zh-cn / de-de status.
# ---------------------------------------------------------------------
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is "b.md".
# ---------------------------------------------------------------------
$zhcn.Range("B3").Value = $status
$zhcn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-02-29 03:55:22"

# ---------------------------------------------------------------------
# de-de sheet: row 3 is "b.md".
# ---------------------------------------------------------------------
$dede.Range("B3").Value = $status
$dede.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("D3").Value = "2016-02-29 03:55:36"

# ---------------------------------------------------------------------
# Hyperlinks: the C3 hyperlink's visible text must match the new handoff
# file name on both sheets. This engine's Hyperlinks collection only
# supports appending/deleting at the sheet level (no true in-place
# update), so each sheet's hyperlinks are rebuilt in original order with
# the single updated display string substituted in for C3.
# ---------------------------------------------------------------------

function Rebuild-Hyperlinks($ws, $links) {
    # Clearing any range's Hyperlinks collection clears the whole sheet,
    # so a single call is enough before re-adding every link.
    $ws.Range("A1").Hyperlinks.Delete()
    foreach ($link in $links) {
        $ws.Hyperlinks.Add($ws.Range($link.Ref), $link.Address, "", "", $link.Display)
    }
}

$zhcnLinks = @(
    @{ Ref = "A2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/6b2327e92806b8cda21b8756b83d118209cd5163/e2e/a.md"; Display = "a.md" },
    @{ Ref = "C2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/78d4b0604ea07b7618673fbe7abd0ad0d9f8688e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Ref = "E2"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3c69a6915739413dc000d50d4435f0fa9b35a0f3/e2e/a.md"; Display = "a.md" },
    @{ Ref = "F2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6e70e1d7a2f2ad56d7af13ff4264f566eccd6a9f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Ref = "A3"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/6b2327e92806b8cda21b8756b83d118209cd5163/e2e/b.md"; Display = "b.md" },
    @{ Ref = "C3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/78d4b0604ea07b7618673fbe7abd0ad0d9f8688e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf" },
    @{ Ref = "E3"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3c69a6915739413dc000d50d4435f0fa9b35a0f3/e2e/a.md"; Display = "a.md" },
    @{ Ref = "F3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6e70e1d7a2f2ad56d7af13ff4264f566eccd6a9f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Ref = "A4"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/6b2327e92806b8cda21b8756b83d118209cd5163/.localization-config"; Display = ".localization-config" }
)

$dedeLinks = @(
    @{ Ref = "A2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/6b2327e92806b8cda21b8756b83d118209cd5163/e2e/a.md"; Display = "a.md" },
    @{ Ref = "C2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1e7ab39da84844b848f0aef1f29b3b1b84daa955/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Ref = "E2"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ef3f82899ce1daa7d7a67364dbeba626c9a67a90/e2e/a.md"; Display = "a.md" },
    @{ Ref = "F2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/690eb82c9c6ec379df585848e59226480018af48/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Ref = "A3"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/6b2327e92806b8cda21b8756b83d118209cd5163/e2e/b.md"; Display = "b.md" },
    @{ Ref = "C3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1e7ab39da84844b848f0aef1f29b3b1b84daa955/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf" },
    @{ Ref = "E3"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ef3f82899ce1daa7d7a67364dbeba626c9a67a90/e2e/a.md"; Display = "a.md" },
    @{ Ref = "F3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/690eb82c9c6ec379df585848e59226480018af48/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Ref = "A4"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/6b2327e92806b8cda21b8756b83d118209cd5163/.localization-config"; Display = ".localization-config" }
)

Rebuild-Hyperlinks $zhcn $zhcnLinks
Rebuild-Hyperlinks $dede $dedeLinks
